$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(71).Insert()

$ws.Cells.Item(71, 1).Value = 4
$ws.Cells.Item(71, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(71, 3).Value = "Los Lagos"
$ws.Cells.Item(71, 4).Value = 44484
$ws.Cells.Item(71, 5).Value = 10
$ws.Cells.Item(71, 6).Value = 100112028
$ws.Cells.Item(71, 7).Value = "Sandia"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 750
$ws.Cells.Item(71, 11).Value = 1000
$ws.Cells.Item(71, 12).Value = 1200
$ws.Cells.Item(71, 13).Value = 1093
$ws.Cells.Item(71, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(71, 15).Value = "Perú"
$ws.Cells.Item(71, 16).Value = 1093
$ws.Cells.Item(71, 17).Value = 1
$ws.Cells.Item(71, 18).Value = "Hortaliza"
